$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.560.11'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '2.505.75'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '316.68'
$ws.Range('E5').Value = '  +3.99%  '
$ws.Range('D6').Value = '95.23'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  -2.04%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('D10').Value = '35.76'
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').Value = '7.53'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').Value = '0.108'
$ws.Range('E13').Value = '  -3.02%  '
$ws.Range('D14').Value = '2.892.59'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '2.524.28'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').Value = '15.18'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').Value = '42.665.56'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').Value = '12.81'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '6.67'
$ws.Range('E20').Value = '  +3.86%  '
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('D22').Value = '69.43'
$ws.Range('E22').Value = '  -2.59%  '
$ws.Range('D23').Value = '249.44'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').Value = '2.95'
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('D25').Value = '2.10'
$ws.Range('E25').Value = '  +3.92%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '26.30'
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = '41.67'
$ws.Range('E28').Value = '  +11.46%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.42'
$ws.Range('E29').Value = '  +4.11%  '
$ws.Range('D30').Value = '10.27'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').Value = '5.93'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = '157.75'
$ws.Range('E32').Value = '  +2.11%  '
$ws.Range('D33').Value = '2.13'
$ws.Range('E33').Value = '  +2.98%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '2.67'
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = '18.91'
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('D36').Value = '3.27'
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('D37').Value = '0.0776'
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('E38').Value = '  -3.52%  '
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('D40').Value = '23.59'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('E41').Value = '  +14.08%  '
$ws.Range('D42').Value = '0.0304'
$ws.Range('E42').Value = '  +1.42%  '
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '3.77'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '3.31'
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('D46').Value = '2.023.98'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('D47').Value = '84.30'
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('D48').Value = '8.89'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('D49').Value = '74.62'
$ws.Range('E49').Value = '  +2.52%  '
$ws.Range('D50').Value = '2.752.35'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.193'
$ws.Range('E51').Value = '  +2.42%  '
